$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 56
$ws.Range("F5").Value = 366
$ws.Range("F6").Value = 173
$ws.Range("F7").Value = 168
$ws.Range("F8").Value = 799
$ws.Range("F9").Value = 4196
$ws.Range("F13").Value = 6055
$ws.Range("F14").Value = 60
$ws.Range("F16").Value = 2330
$ws.Range("F18").Value = 166
$ws.Range("F19").Value = 470
$ws.Range("F20").Value = 9118
$ws.Range("F22").Value = 2452
$ws.Range("F24").Value = 2308
$ws.Range("F25").Value = 2432
$ws.Range("F26").Value = 1390
$ws.Range("F27").Value = 239
$ws.Range("F28").Value = 1957
$ws.Range("F31").Value = 325
$ws.Range("F33").Value = 37
$ws.Range("F34").Value = 279
$ws.Range("F35").Value = 42
$ws.Range("F36").Value = 50
$ws.Range("F37").Value = 383
$ws.Range("F38").Value = 1221
$ws.Range("F40").Value = 73
$ws.Range("F41").Value = 96
$ws.Range("F42").Value = 236
$ws.Range("F43").Value = 1529
$ws.Range("F44").Value = 2499
$ws.Range("F45").Value = 921
$ws.Range("F48").Value = 16

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 168
$ws.Range("F22").Value = 65
$ws.Range("F23").Value = 65

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 690
$ws.Range("F3").Value = 888
$ws.Range("F4").Value = 101

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 56
$ws.Range("F4").Value = 690
$ws.Range("F5").Value = 888
$ws.Range("F6").Value = 101
$ws.Range("F7").Value = 366
$ws.Range("F9").Value = 173
$ws.Range("F11").Value = 168
$ws.Range("F12").Value = 168
$ws.Range("F13").Value = 799
$ws.Range("F14").Value = 4196
$ws.Range("F16").Value = 6055
$ws.Range("F17").Value = 60
$ws.Range("F19").Value = 2330
$ws.Range("F20").Value = 166
$ws.Range("F21").Value = 470
$ws.Range("F22").Value = 9118
$ws.Range("F24").Value = 2452
$ws.Range("F25").Value = 2308
$ws.Range("F26").Value = 1390
$ws.Range("F27").Value = 239
$ws.Range("F28").Value = 1957
$ws.Range("F31").Value = 325
$ws.Range("F32").Value = 37
$ws.Range("F33").Value = 279
$ws.Range("F34").Value = 42
$ws.Range("F35").Value = 50
$ws.Range("F36").Value = 383
$ws.Range("F38").Value = 73
$ws.Range("F39").Value = 96
$ws.Range("F40").Value = 236
$ws.Range("F41").Value = 1529
$ws.Range("F42").Value = 2499
$ws.Range("F43").Value = 921
$ws.Range("F49").Value = 16
$ws.Range("F50").Value = 65
